# Auto-generated edit script
# Updates market-price-derived Leve profit columns (H-N) for various rows
# across all 8 worksheets, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132: Fast-forwarding Flora
$ws.Cells.Item(132, 8).Value = 1566.9122
$ws.Cells.Item(132, 9).Value = 1617.9215
$ws.Cells.Item(132, 10).Value = 1133.3334
$ws.Cells.Item(132, 11).Value = 4853.764499999999
$ws.Cells.Item(132, 12).Value = 3400.0002
$ws.Cells.Item(132, 13).Value = -2323.764499999999
$ws.Cells.Item(132, 14).Value = -8460.0002

# Row 134: Binding Spells
$ws.Cells.Item(134, 8).Value = 85618.17999999999
$ws.Cells.Item(134, 10).Value = 85618.17999999999
$ws.Cells.Item(134, 12).Value = 85618.17999999999
$ws.Cells.Item(134, 14).Value = -95758.17999999999

# Row 135: For Tired Minds
$ws.Cells.Item(135, 8).Value = 736.4386
$ws.Cells.Item(135, 9).Value = 388.375
$ws.Cells.Item(135, 10).Value = 2592.7778
$ws.Cells.Item(135, 11).Value = 3495.375
$ws.Cells.Item(135, 12).Value = 23335.0002
$ws.Cells.Item(135, 13).Value = -960.375
$ws.Cells.Item(135, 14).Value = -28405.0002

# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 810.29034
$ws.Cells.Item(137, 9).Value = 732.4211
$ws.Cells.Item(137, 10).Value = 933.5833
$ws.Cells.Item(137, 11).Value = 2197.2633
$ws.Cells.Item(137, 12).Value = 2800.7499
$ws.Cells.Item(137, 13).Value = 352.7366999999999
$ws.Cells.Item(137, 14).Value = -7900.7499

# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 1455.49
$ws.Cells.Item(138, 9).Value = 765.2558
$ws.Cells.Item(138, 10).Value = 1976.193
$ws.Cells.Item(138, 11).Value = 2295.7674
$ws.Cells.Item(138, 12).Value = 5928.579
$ws.Cells.Item(138, 13).Value = 2844.2326
$ws.Cells.Item(138, 14).Value = -16208.579

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 9707.17
$ws.Cells.Item(32, 9).Value = 6195.6978
$ws.Cells.Item(32, 10).Value = 31277.643
$ws.Cells.Item(32, 11).Value = 6195.6978
$ws.Cells.Item(32, 12).Value = 31277.643
$ws.Cells.Item(32, 13).Value = -5908.6978
$ws.Cells.Item(32, 14).Value = -31851.643

# Row 45: Hollow Hallmarks
$ws.Cells.Item(45, 8).Value = 1644.7878
$ws.Cells.Item(45, 9).Value = 927.9524
$ws.Cells.Item(45, 10).Value = 2899.25
$ws.Cells.Item(45, 11).Value = 927.9524
$ws.Cells.Item(45, 12).Value = 2899.25
$ws.Cells.Item(45, 13).Value = -550.9524
$ws.Cells.Item(45, 14).Value = -3653.25

# Row 132: Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 2505.318
$ws.Cells.Item(132, 9).Value = 2305.9
$ws.Cells.Item(132, 10).Value = 4499.5
$ws.Cells.Item(132, 11).Value = 6917.700000000001
$ws.Cells.Item(132, 12).Value = 13498.5
$ws.Cells.Item(132, 13).Value = -4387.700000000001
$ws.Cells.Item(132, 14).Value = -18558.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Cells.Item(20, 8).Value = 26396.27
$ws.Cells.Item(20, 9).Value = 60188.727
$ws.Cells.Item(20, 10).Value = 1615.1333
$ws.Cells.Item(20, 11).Value = 60188.727
$ws.Cells.Item(20, 12).Value = 1615.1333
$ws.Cells.Item(20, 13).Value = -59941.727
$ws.Cells.Item(20, 14).Value = -2109.1333

# Row 134: Ruthenium Supremium
$ws.Cells.Item(134, 8).Value = 26810.244
$ws.Cells.Item(134, 9).Value = 2098.818
$ws.Cells.Item(134, 10).Value = 128744.875
$ws.Cells.Item(134, 11).Value = 6296.454000000001
$ws.Cells.Item(134, 12).Value = 386234.625
$ws.Cells.Item(134, 13).Value = -3761.454000000001
$ws.Cells.Item(134, 14).Value = -391304.625

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Cells.Item(16, 8).Value = 1710.1818
$ws.Cells.Item(16, 9).Value = 850
$ws.Cells.Item(16, 10).Value = 4004
$ws.Cells.Item(16, 11).Value = 850
$ws.Cells.Item(16, 12).Value = 4004
$ws.Cells.Item(16, 13).Value = -563
$ws.Cells.Item(16, 14).Value = -4578

# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 1733014.8
$ws.Cells.Item(31, 9).Value = 2144941
$ws.Cells.Item(31, 11).Value = 2144941
$ws.Cells.Item(31, 13).Value = -2144646

# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 1733014.8
$ws.Cells.Item(34, 9).Value = 2144941
$ws.Cells.Item(34, 11).Value = 2144941
$ws.Cells.Item(34, 13).Value = -2144739

# Row 58: You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 3722.675
$ws.Cells.Item(58, 9).Value = 1221.5
$ws.Cells.Item(58, 11).Value = 1221.5
$ws.Cells.Item(58, 13).Value = -1018.5

# Row 113: Patient Patients
$ws.Cells.Item(113, 8).Value = 1710.1818
$ws.Cells.Item(113, 9).Value = 850
$ws.Cells.Item(113, 10).Value = 4004
$ws.Cells.Item(113, 11).Value = 850
$ws.Cells.Item(113, 12).Value = 4004
$ws.Cells.Item(113, 13).Value = 1320
$ws.Cells.Item(113, 14).Value = -8344

# Row 122: Timber of Tenkonto
$ws.Cells.Item(122, 8).Value = 961.0526
$ws.Cells.Item(122, 9).Value = 604.6667
$ws.Cells.Item(122, 10).Value = 1572
$ws.Cells.Item(122, 11).Value = 1814.0001
$ws.Cells.Item(122, 12).Value = 4716
$ws.Cells.Item(122, 13).Value = 635.9999
$ws.Cells.Item(122, 14).Value = -9616

# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 1491.025
$ws.Cells.Item(132, 9).Value = 916.42554
$ws.Cells.Item(132, 10).Value = 2309.394
$ws.Cells.Item(132, 11).Value = 2749.27662
$ws.Cells.Item(132, 12).Value = 6928.181999999999
$ws.Cells.Item(132, 13).Value = -219.2766199999996
$ws.Cells.Item(132, 14).Value = -11988.182

# Row 134: Wood You Be Quiet
$ws.Cells.Item(134, 8).Value = 1446.68
$ws.Cells.Item(134, 9).Value = 1381.2222
$ws.Cells.Item(134, 10).Value = 1615
$ws.Cells.Item(134, 11).Value = 4143.6666
$ws.Cells.Item(134, 12).Value = 4845
$ws.Cells.Item(134, 13).Value = -1608.6666
$ws.Cells.Item(134, 14).Value = -9915

# Row 136: Turali Quality
$ws.Cells.Item(136, 8).Value = 3722.675
$ws.Cells.Item(136, 9).Value = 1221.5
$ws.Cells.Item(136, 11).Value = 3664.5
$ws.Cells.Item(136, 13).Value = -1114.5

$ws = $wb.Worksheets.Item("CUL")
# Row 14: Keep Your Powder Dry
$ws.Cells.Item(14, 8).Value = 144.13333
$ws.Cells.Item(14, 9).Value = 144.13333
$ws.Cells.Item(14, 11).Value = 432.39999
$ws.Cells.Item(14, 13).Value = -259.39999

# Row 117: A Good Omen
$ws.Cells.Item(117, 8).Value = 2539.2942
$ws.Cells.Item(117, 9).Value = 1500
$ws.Cells.Item(117, 10).Value = 2604.25
$ws.Cells.Item(117, 11).Value = 4500
$ws.Cells.Item(117, 12).Value = 7812.75
$ws.Cells.Item(117, 13).Value = -1058
$ws.Cells.Item(117, 14).Value = -14696.75

# Row 132: More Mezcal
$ws.Cells.Item(132, 8).Value = 762
$ws.Cells.Item(132, 9).Value = 611.0769
$ws.Cells.Item(132, 10).Value = 1252.5
$ws.Cells.Item(132, 11).Value = 5499.6921
$ws.Cells.Item(132, 12).Value = 11272.5
$ws.Cells.Item(132, 13).Value = -2969.6921
$ws.Cells.Item(132, 14).Value = -16332.5

# Row 140: Sweet, Sweet Bean Juice
$ws.Cells.Item(140, 8).Value = 113655.26
$ws.Cells.Item(140, 9).Value = 144751.58
$ws.Cells.Item(140, 10).Value = 4818.1665
$ws.Cells.Item(140, 11).Value = 434254.74
$ws.Cells.Item(140, 12).Value = 14454.4995
$ws.Cells.Item(140, 13).Value = -429074.74
$ws.Cells.Item(140, 14).Value = -24814.4995

$ws = $wb.Worksheets.Item("GSM")
# Row 11: A Ringing Success
$ws.Cells.Item(11, 8).Value = 15613200
$ws.Cells.Item(11, 9).Value = 18262750
$ws.Cells.Item(11, 10).Value = 5015000
$ws.Cells.Item(11, 11).Value = 18262750
$ws.Cells.Item(11, 12).Value = 5015000
$ws.Cells.Item(11, 13).Value = -18262611
$ws.Cells.Item(11, 14).Value = -5015278

# Row 127: Sage with the Golden Earrings
$ws.Cells.Item(127, 8).Value = 31660
$ws.Cells.Item(127, 10).Value = 31660
$ws.Cells.Item(127, 12).Value = 31660
$ws.Cells.Item(127, 14).Value = -41580

# Row 132: On Board for Lar
$ws.Cells.Item(132, 8).Value = 2322.025
$ws.Cells.Item(132, 9).Value = 2188.7827
$ws.Cells.Item(132, 10).Value = 2502.2942
$ws.Cells.Item(132, 11).Value = 6566.348100000001
$ws.Cells.Item(132, 12).Value = 7506.882599999999
$ws.Cells.Item(132, 13).Value = -4036.348100000001
$ws.Cells.Item(132, 14).Value = -12566.8826

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 2315.9302
$ws.Cells.Item(132, 9).Value = 2114.7273
$ws.Cells.Item(132, 10).Value = 2979.9
$ws.Cells.Item(132, 11).Value = 6344.1819
$ws.Cells.Item(132, 12).Value = 8939.700000000001
$ws.Cells.Item(132, 13).Value = -3814.1819
$ws.Cells.Item(132, 14).Value = -13999.7

# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 2416.4
$ws.Cells.Item(136, 9).Value = 1250.7742
$ws.Cells.Item(136, 10).Value = 11450
$ws.Cells.Item(136, 11).Value = 3752.3226
$ws.Cells.Item(136, 12).Value = 34350
$ws.Cells.Item(136, 13).Value = -1202.3226
$ws.Cells.Item(136, 14).Value = -39450

$ws = $wb.Worksheets.Item("WVR")
# Row 110: Suits You
$ws.Cells.Item(110, 8).Value = 28500
$ws.Cells.Item(110, 10).Value = 28500
$ws.Cells.Item(110, 12).Value = 28500
$ws.Cells.Item(110, 14).Value = -36680

# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 1083.9615
$ws.Cells.Item(132, 9).Value = 888.7
$ws.Cells.Item(132, 10).Value = 1734.8334
$ws.Cells.Item(132, 11).Value = 2666.1
$ws.Cells.Item(132, 12).Value = 5204.5002
$ws.Cells.Item(132, 13).Value = -136.1000000000004
$ws.Cells.Item(132, 14).Value = -10264.5002

# Row 136: Weaving the Envelope
$ws.Cells.Item(136, 8).Value = 1248.6154
$ws.Cells.Item(136, 9).Value = 1447.8422
$ws.Cells.Item(136, 10).Value = 707.8570999999999
$ws.Cells.Item(136, 11).Value = 4343.5266
$ws.Cells.Item(136, 12).Value = 2123.5713
$ws.Cells.Item(136, 13).Value = -1793.5266
$ws.Cells.Item(136, 14).Value = -7223.5713
